# B6-PowerPoint.pptx edit
# 1) Re-style the three tables (slides 14-16) with the new built-in table style.
# 2) Re-colour the deck's (slide-master) theme colour scheme from the
#    "Red Violet"/Integral palette to the standard "Office" palette -
#    this is the reachable, colour-scheme part of the theme1.xml/theme2.xml
#    swap described by the diff (the Notes Master's theme, theme1.xml, is
#    not reachable from the exposed PowerPoint object model).

$p = $ppt.ActivePresentation

$newTableStyleId = "{56477778-ABE0-441D-B949-2EB91E7FA8A8}"

foreach ($slideIdx in 14,15,16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# Office theme colours (RRGGBB) keyed by ThemeColorScheme.Colors() index.
# COM RGB values are COLORREF (0x00BBGGRR), so each entry below is built
# from the R/G/B bytes of the target hex colour rather than the hex
# literal itself.
function New-ComRgb([int]$r, [int]$g, [int]$b) {
    return $r -bor ($g * 256) -bor ($b * 65536)
}

$officeColors = @{
    1  = (New-ComRgb 0x00 0x00 0x00)  # dk1
    2  = (New-ComRgb 0xFF 0xFF 0xFF)  # lt1
    3  = (New-ComRgb 0x44 0x54 0x6A)  # dk2
    4  = (New-ComRgb 0xE7 0xE6 0xE6)  # lt2
    5  = (New-ComRgb 0x5B 0x9B 0xD5)  # accent1
    6  = (New-ComRgb 0xED 0x7D 0x31)  # accent2
    7  = (New-ComRgb 0xA5 0xA5 0xA5)  # accent3
    8  = (New-ComRgb 0xFF 0xC0 0x00)  # accent4
    9  = (New-ComRgb 0x44 0x72 0xC4)  # accent5
    10 = (New-ComRgb 0x70 0xAD 0x47)  # accent6
    11 = (New-ComRgb 0x05 0x63 0xC1)  # hlink
    12 = (New-ComRgb 0x95 0x4F 0x72)  # folHlink
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i]
}
